$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation for "Vega Modelo de Temuco - Betarraga" needs to be
# inserted into the daily log at row 167 (dated 2021-09-16 / serial 44455),
# pushing the existing rows 167-248 down to 168-249.
$ws.Rows.Item(167).Insert()

# The freshly inserted row 167 is blank; seed it with the same reference data
# (market/category/quality/price-unit/origin columns) as the row immediately
# below it (which now holds what used to be the old row 167), then overwrite
# just the date (D) and volume (J) with the new observation's values.
$ws.Range("A168:R168").Copy()
$ws.Range("A167").PasteSpecial()

$ws.Range("D167").Value = 44455
$ws.Range("J167").Value = 20
